$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that changed, per the repull/recalculation.
$ws.Range("F2").Value = -3
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = -1
$ws.Range("F10").Value = -3
$ws.Range("F12").Value = 6
